$wb = $excel.ActiveWorkbook

# Common message-sheet values (reused across CypherOutput_Message and StatOutput_Message)
$neo4jUrlLabel = 'Neo4j_URL:'
$neo4jUrlValue = 'bolt://ncias-q2251-c.nci.nih.gov:7687'
$userNameLabel = 'User_name:'
$userNameValue = 'neo4j'
$pwdLabel = 'PWD:'
$pwdValue = 'icdcDBneo4j0'
$cypherLabel = 'Cypher:'
$cypher1 = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN [''Malignant neoplasm of the respiratory tract cell type specified :: Lung adenocarcinoma (metastatic)''] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(s.clinical_study_designation,'''') AS `Study Code` , coalesce(s.clinical_study_type,'''') AS  `Study Type`, coalesce(demo.breed,'''') AS Breed , coalesce(diag.disease_term,'''') AS Diagnosis , coalesce(diag.stage_of_disease,'''') AS `Stage of Disease` ,  coalesce(demo.patient_age_at_enrollment,'''') AS Age , coalesce(demo.sex,'''') AS Sex , coalesce(demo.neutered_indicator,'''') AS  `Neutered Status`'
$outputLabel = 'Output:'
$outputPath = 'C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC12_Canine_Filter_Diagnosis-MaligMetastatic_Neo4jData.xlsx'
$cypher2 = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN[''Malignant neoplasm of the respiratory tract cell type specified :: Lung adenocarcinoma (metastatic)'']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study'

# ---- Sheet: CypherOutput_Message (copy of Message sheet) ----
$wsCypherMsg = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsCypherMsg.Name = "CypherOutput_Message"
$wsCypherMsg.Range("A1").Value = $neo4jUrlLabel
$wsCypherMsg.Range("A2").Value = $neo4jUrlValue
$wsCypherMsg.Range("A3").Value = $userNameLabel
$wsCypherMsg.Range("A4").Value = $userNameValue
$wsCypherMsg.Range("A5").Value = $pwdLabel
$wsCypherMsg.Range("A6").Value = $pwdValue
$wsCypherMsg.Range("A7").Value = $cypherLabel
$wsCypherMsg.Range("A8").Value = $cypher1
$wsCypherMsg.Range("A9").Value = $outputLabel
$wsCypherMsg.Range("A10").Value = $outputPath

# ---- Sheet: StatOutput ----
$wsStatOutput = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsStatOutput.Name = "StatOutput"
$wsStatOutput.Range("A1").Value = "number_of_files"
$wsStatOutput.Range("B1").Value = "number_of_sample"
$wsStatOutput.Range("C1").Value = "number_of_cases"
$wsStatOutput.Range("D1").Value = "number_of_study"
$wsStatOutput.Range("A2").Formula = "'3"
$wsStatOutput.Range("B2").Formula = "'5"
$wsStatOutput.Range("C2").Formula = "'2"
$wsStatOutput.Range("D2").Formula = "'1"

# ---- Sheet: StatOutput_Message (Message sheet content twice, 2nd pass uses new Cypher) ----
$wsStatMsg = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsStatMsg.Name = "StatOutput_Message"
$wsStatMsg.Range("A1").Value = $neo4jUrlLabel
$wsStatMsg.Range("A2").Value = $neo4jUrlValue
$wsStatMsg.Range("A3").Value = $userNameLabel
$wsStatMsg.Range("A4").Value = $userNameValue
$wsStatMsg.Range("A5").Value = $pwdLabel
$wsStatMsg.Range("A6").Value = $pwdValue
$wsStatMsg.Range("A7").Value = $cypherLabel
$wsStatMsg.Range("A8").Value = $cypher1
$wsStatMsg.Range("A9").Value = $outputLabel
$wsStatMsg.Range("A10").Value = $outputPath
$wsStatMsg.Range("A11").Value = $neo4jUrlLabel
$wsStatMsg.Range("A12").Value = $neo4jUrlValue
$wsStatMsg.Range("A13").Value = $userNameLabel
$wsStatMsg.Range("A14").Value = $userNameValue
$wsStatMsg.Range("A15").Value = $pwdLabel
$wsStatMsg.Range("A16").Value = $pwdValue
$wsStatMsg.Range("A17").Value = $cypherLabel
$wsStatMsg.Range("A18").Value = $cypher2
$wsStatMsg.Range("A19").Value = $outputLabel
$wsStatMsg.Range("A20").Value = $outputPath

# Restore the originally-active sheet/tab (CypherOutput) so the workbook-level
# view state matches the pre-edit workbook (only new sheets were appended).
$wb.Worksheets.Item(1).Activate()
